# "revert arguments parsing changes"
#
# The "hasTotalSteps" step (row 21) and its paired log row (row 22) in
# the "Case" sheet recorded the `numberOfSteps` argument as a bare
# number (numberOfSteps:0). This reverts that to the original quoted
# string-literal form (numberOfSteps:'0') used everywhere else in the
# sheet's arguments parsing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D21").Value = "object=Test;name=test-1;action=hasTotalSteps;arguments={numberOfSteps:'0'}"
$ws.Range("D22").Value = "row=6;arguments={numberOfSteps:'0'}"

# Move the active selection to D22 to match the saved view state.
$ws.Range("D22").Select()
